$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, image (B), word (C), category (D)
$data = @(
    @(2, "face/face020.jpg", "tagen", "face"),
    @(3, "face/face017.jpg", "schätzen", "face"),
    @(4, "flower/flower000.jpg", "drohen", "flower"),
    @(5, "face/face030.jpg", "klappen", "face"),
    @(6, "face/face005.jpg", "posten", "face"),
    @(7, "face/face022.jpg", "fliegen", "face"),
    @(8, "flower/flower017.jpg", "nehmen", "flower"),
    @(9, "flower/flower030.jpg", "hauen", "flower"),
    @(10, "flower/flower026.jpg", "fesseln", "flower"),
    @(11, "flower/flower023.jpg", "scheitern", "flower"),
    @(12, "flower/flower010.jpg", "fliehen", "flower"),
    @(13, "flower/flower018.jpg", "sondern", "flower"),
    @(14, "flower/flower001.jpg", "saufen", "flower"),
    @(15, "face/face003.jpg", "biegen", "face"),
    @(16, "face/face011.jpg", "ehren", "face"),
    @(17, "face/face000.jpg", "backen", "face"),
    @(18, "face/face029.jpg", "hupen", "face"),
    @(19, "face/face015.jpg", "schenken", "face"),
    @(20, "flower/flower024.jpg", "wenden", "flower"),
    @(21, "face/face013.jpg", "starten", "face"),
    @(22, "face/face004.jpg", "husten", "face"),
    @(23, "flower/flower009.jpg", "füttern", "flower"),
    @(24, "face/face021.jpg", "antun", "face"),
    @(25, "face/face027.jpg", "segeln", "face"),
    @(26, "face/face012.jpg", "strahlen", "face"),
    @(27, "flower/flower007.jpg", "regnen", "flower"),
    @(28, "flower/flower013.jpg", "rücken", "flower"),
    @(29, "flower/flower004.jpg", "dauern", "flower"),
    @(30, "flower/flower002.jpg", "mieten", "flower"),
    @(31, "flower/flower012.jpg", "langen", "flower"),
    @(32, "face/face028.jpg", "jubeln", "face"),
    @(33, "flower/flower028.jpg", "fühlen", "flower"),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

Write-Output "done"